# "merge fail - copy in files" — the NBACTR (Non-BAU Average Corporate Tax
# Rate) sheet was pulled into this workbook by mistake during a bad merge.
# Back it out: drop the extra sheet, remove the "About" blurb row that
# pointed at it, and let the shared-string table shrink accordingly.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The "About" sheet's row 2 held the NBACTR description line; deleting the
# whole row shifts everything below it up by one (row 4 -> 3, ... row 11 -> 10).
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows(2).Delete()

# Remove the NBACTR worksheet entirely - it shouldn't have been copied in.
$wsNBACTR = $wb.Worksheets.Item("NBACTR")
[void]$wsNBACTR.Delete()

# Keep "About" the active sheet/selection, now pointed at its last row
# (previously A3/A11, shifted up along with the row delete above).
$wsAbout.Activate()
$wsAbout.Range("A11").Select() | Out-Null
